$d = $word.ActiveDocument

# Delete the first 5 paragraphs (the citation list entries), leaving only
# the trailing empty "TS"-styled paragraph before the section break.
for ($i = 1; $i -le 5; $i++) {
    $d.Paragraphs(1).Range.Delete()
}
